$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the values of columns B, D, E, G, H between row 2 and row 3
$cols = @("B", "D", "E", "G", "H")

foreach ($col in $cols) {
    $addr2 = $col + "2"
    $addr3 = $col + "3"
    $val2 = $ws.Range($addr2).Value()
    $val3 = $ws.Range($addr3).Value()
    $ws.Range($addr2).Value = $val3
    $ws.Range($addr3).Value = $val2
}

$wb.Save()
